$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1497.4375
$ws.Range("I19").Value = 1498.4
$ws.Range("J19").Value = 1483
$ws.Range("K19").Value = 1498.4
$ws.Range("L19").Value = 1483
$ws.Range("M19").Value = -1323.4
$ws.Range("N19").Value = -1833

$ws.Range("H64").Value = 6832.2666
$ws.Range("I64").Value = 4580.5
$ws.Range("K64").Value = 4580.5
$ws.Range("M64").Value = -4332.5

$ws.Range("H67").Value = 6832.2666
$ws.Range("I67").Value = 4580.5
$ws.Range("K67").Value = 4580.5
$ws.Range("M67").Value = -3722.5

$ws.Range("H116").Value = 21750906
$ws.Range("I116").Value = 29425616
$ws.Range("K116").Value = 29425616
$ws.Range("M116").Value = -29422174

$ws.Range("H129").Value = 1317.9445
$ws.Range("I129").Value = 772.3
$ws.Range("K129").Value = 2316.9
$ws.Range("M129").Value = 2683.1

$ws.Range("H132").Value = 1862.2115
$ws.Range("I132").Value = 1870.4131
$ws.Range("K132").Value = 5611.2393
$ws.Range("M132").Value = -3081.2393

$ws.Range("H138").Value = 1978.31
$ws.Range("I138").Value = 706.4865
$ws.Range("J138").Value = 2725.254
$ws.Range("K138").Value = 2119.4595
$ws.Range("L138").Value = 8175.762
$ws.Range("M138").Value = 3020.5405
$ws.Range("N138").Value = -18455.762

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2853.0488
$ws.Range("I61").Value = 2406.8845
$ws.Range("J61").Value = 3626.4
$ws.Range("K61").Value = 2406.8845
$ws.Range("L61").Value = 3626.4
$ws.Range("M61").Value = -2194.8845
$ws.Range("N61").Value = -4050.4

$ws.Range("H63").Value = 4484.5713
$ws.Range("I63").Value = 2000
$ws.Range("J63").Value = 5162.1816
$ws.Range("K63").Value = 2000
$ws.Range("L63").Value = 5162.1816
$ws.Range("M63").Value = -1314
$ws.Range("N63").Value = -6534.1816

$ws.Range("H66").Value = 4484.5713
$ws.Range("I66").Value = 2000
$ws.Range("J66").Value = 5162.1816
$ws.Range("K66").Value = 10000
$ws.Range("L66").Value = 25810.908
$ws.Range("M66").Value = -6568
$ws.Range("N66").Value = -32674.908

$ws.Range("H74").Value = 2753.9429
$ws.Range("I74").Value = 2243
$ws.Range("J74").Value = 5819.6
$ws.Range("K74").Value = 2243
$ws.Range("L74").Value = 5819.6
$ws.Range("M74").Value = -1369
$ws.Range("N74").Value = -7567.6

$ws.Range("H77").Value = 2753.9429
$ws.Range("I77").Value = 2243
$ws.Range("J77").Value = 5819.6
$ws.Range("K77").Value = 11215
$ws.Range("L77").Value = 29098
$ws.Range("M77").Value = -6847
$ws.Range("N77").Value = -37834

$ws.Range("H136").Value = 2853.0488
$ws.Range("I136").Value = 2406.8845
$ws.Range("J136").Value = 3626.4
$ws.Range("K136").Value = 7220.6535
$ws.Range("L136").Value = 10879.2
$ws.Range("M136").Value = -4670.6535
$ws.Range("N136").Value = -15979.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2508.111
$ws.Range("I99").Value = 1564.7
$ws.Range("K99").Value = 1564.7
$ws.Range("M99").Value = -66.70000000000005

$ws.Range("H134").Value = 1883002.2
$ws.Range("I134").Value = 2465134.2
$ws.Range("K134").Value = 7395402.600000001
$ws.Range("M134").Value = -7392867.600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4991.5312
$ws.Range("I31").Value = 1990.2632
$ws.Range("J31").Value = 9378
$ws.Range("K31").Value = 1990.2632
$ws.Range("L31").Value = 9378
$ws.Range("M31").Value = -1695.2632
$ws.Range("N31").Value = -9968

$ws.Range("H34").Value = 4991.5312
$ws.Range("I34").Value = 1990.2632
$ws.Range("J34").Value = 9378
$ws.Range("K34").Value = 1990.2632
$ws.Range("L34").Value = 9378
$ws.Range("M34").Value = -1788.2632
$ws.Range("N34").Value = -9782

$ws.Range("H132").Value = 2986.85
$ws.Range("I132").Value = 2929.375
$ws.Range("K132").Value = 8788.125
$ws.Range("M132").Value = -6258.125

$ws.Range("H134").Value = 3732.1667
$ws.Range("J134").Value = 4200
$ws.Range("L134").Value = 12600
$ws.Range("N134").Value = -17670

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 6738
$ws.Range("I122").Value = 12984.333
$ws.Range("J122").Value = 491.66666
$ws.Range("K122").Value = 116858.997
$ws.Range("L122").Value = 4424.99994
$ws.Range("M122").Value = -114408.997
$ws.Range("N122").Value = -9324.99994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3999
$ws.Range("I80").Value = 3999
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3999
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3001
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 3999
$ws.Range("I83").Value = 3999
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 19995
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -15003
$ws.Range("N83").ClearContents()

$ws.Range("H97").Value = 1127.875
$ws.Range("J97").Value = 1537
$ws.Range("L97").Value = 1537
$ws.Range("N97").Value = -2529

$ws.Range("H99").Value = 32293
$ws.Range("I99").Value = 13866.25
$ws.Range("J99").Value = 106000
$ws.Range("K99").Value = 13866.25
$ws.Range("L99").Value = 106000
$ws.Range("M99").Value = -11620.25
$ws.Range("N99").Value = -110492

$ws.Range("H132").Value = 4076.611
$ws.Range("I132").Value = 3348.7144
$ws.Range("K132").Value = 10046.1432
$ws.Range("M132").Value = -7516.143199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 22226240
$ws.Range("I40").Value = 33336162
$ws.Range("K40").Value = 33336162
$ws.Range("M40").Value = -33336026

$ws.Range("H82").Value = 4044.7856
$ws.Range("I82").Value = 4144
$ws.Range("J82").Value = 3866.2
$ws.Range("K82").Value = 4144
$ws.Range("L82").Value = 3866.2
$ws.Range("M82").Value = -3783
$ws.Range("N82").Value = -4588.2

$ws.Range("H85").Value = 4044.7856
$ws.Range("I85").Value = 4144
$ws.Range("J85").Value = 3866.2
$ws.Range("K85").Value = 4144
$ws.Range("L85").Value = 3866.2
$ws.Range("M85").Value = -2896
$ws.Range("N85").Value = -6362.2

$ws.Range("H100").Value = 2810.5
$ws.Range("I100").Value = 2810.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2810.5
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2269.5
$ws.Range("N100").ClearContents()

$ws.Range("H132").Value = 3242.0625
$ws.Range("I132").Value = 3124.8667
$ws.Range("K132").Value = 9374.6001
$ws.Range("M132").Value = -6844.6001

$ws.Range("H136").Value = 1929.2354
$ws.Range("I136").Value = 1719.8
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 5159.4
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -2609.4
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 11135
$ws.Range("J32").Value = 10980
$ws.Range("L32").Value = 10980
$ws.Range("N32").Value = -11614

$ws.Range("H96").Value = 9547
$ws.Range("I96").Value = 4924.5
$ws.Range("J96").Value = 14169.5
$ws.Range("K96").Value = 4924.5
$ws.Range("L96").Value = 14169.5
$ws.Range("M96").Value = -3551.5
$ws.Range("N96").Value = -16915.5
